$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.018.86'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '2.309.65'
$ws.Range('E3').Value = '  -3.06%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'303.51"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.20%  '
$ws.Range('D6').Value = "'99.66"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.41%  '
$ws.Range('D7').Value = "'0.509"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.72%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('D10').Value = "'35.09"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.24%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = "'0.0797"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.21%  '
$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D12').Value = "'51.21"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.25%  '
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').Value = "'6.79"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.29%  '
$ws.Range('D15').Value = '2.668.30'
$ws.Range('E15').Value = '  -3.05%  '
$ws.Range('D16').Value = "'15.55"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.81%  '
$ws.Range('D17').Value = '2.297.57'
$ws.Range('E17').Value = '  -3.51%  '
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('D19').Value = '42.960.92'
$ws.Range('E19').Value = '  -1.48%  '
$ws.Range('D20').Value = "'11.73"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.35%  '
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('D22').Value = "'6.06"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.94%  '
$ws.Range('D23').Value = "'67.53"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('D24').Value = "'236.77"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('D25').Value = "'1.97"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.71%  '
$ws.Range('D26').Value = "'2.53"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.93%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = "'24.94"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.55%  '
$ws.Range('D29').Value = "'2.18"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.70%  '
$ws.Range('D30').Value = "'34.74"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.56%  '
$ws.Range('D31').Value = "'165.10"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.94%  '
$ws.Range('D32').Value = "'9.18"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.19%  '
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').Value = "'5.04"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.63%  '
$ws.Range('D35').Value = "'2.41"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.85%  '
$ws.Range('D36').Value = "'4.48"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.47%  '
$ws.Range('E37').Value = '  -4.88%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = "'2.87"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.47%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = "'16.46"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -10.79%  '
$ws.Range('E40').Value = '  -7.23%  '
$ws.Range('D41').Value = "'0.101"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.20%  '
$ws.Range('E42').Value = '  -2.69%  '
$ws.Range('D43').Value = "'2.43"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.86%  '
$ws.Range('D44').Value = '1.977.98'
$ws.Range('E44').Value = '  -2.93%  '
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('D46').Value = "'18.57"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Value = "'9.88"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.70%  '
$ws.Range('D48').Value = "'2.90"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -8.24%  '
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('D50').Value = "'53.97"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.15%  '
$ws.Range('D51').Value = '2.536.23'
$ws.Range('E51').Value = '  -2.52%  '
